# This script rearranges the species-observation records found in rows 4-8
# of the active worksheet. Only columns A, B, D, E, F, G, H, Q, R participate
# in the rearrangement; all other columns (C, I, K, P, S:W, Y:AB, AD, AE, AG,
# AT, AW:AY) are identical across these rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are being permuted between rows 4-8.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Capture the current ("before") values for each row/column pair first,
# since the rearrangement below would otherwise clobber source data before
# it has been read.
$snapshot = @{}
foreach ($r in 4..8) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Target row <- source row mapping (i.e. row $target ends up holding the
# values that currently live in row $source).
$mapping = @{
    4 = 5
    5 = 8
    6 = 7
    7 = 4
    8 = 6
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $srcVals = $snapshot[$source]
    foreach ($col in $cols) {
        $ws.Range("$col$target").Value = $srcVals[$col]
    }
}
